# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2 = "H" home splits)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 394
$wsOff.Range("C2").Value = 277
$wsOff.Range("D2").Value = 99
$wsOff.Range("E2").Value = 45
$wsOff.Range("F2").Value = 10
$wsOff.Range("G2").Value = 7

# Update DEF sheet (row 2 = "H" home splits)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 518
$wsDef.Range("C2").Value = 359
$wsDef.Range("D2").Value = 133
$wsDef.Range("E2").Value = 57
$wsDef.Range("F2").Value = 5
$wsDef.Range("G2").Value = 6
